# Natmi following Dr Hou advice
# Update C1qa-Cspg4 LR-pair sheet to include all sending/target cluster combinations
# (ECs and M2 as sending clusters; ECs, FAPs, M2, sCs as target clusters)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 8 data rows (rows 2-9) x 20 columns (A-T)
$data = New-Object 'object[,]' 8,20

# Row 2: ECs -> ECs
$data[0,0] = "ECs"   # Sending cluster
$data[0,1] = "C1qa"   # Ligand symbol
$data[0,2] = "Cspg4"   # Receptor symbol
$data[0,3] = "ECs"   # Target cluster
$data[0,4] = 1   # Ligand-expressing cells
$data[0,5] = 0.3333333333333333   # Ligand detection rate
$data[0,6] = 5.392385333333333   # Ligand average expression value
$data[0,7] = 16.177156   # Ligand total expression value
$data[0,8] = 0.04998147672264548   # Ligand derived specificity of average expression value
$data[0,9] = 0.04998147672264548   # Ligand derived specificity of total expression value
$data[0,10] = 3   # Receptor-expressing cells
$data[0,11] = 1   # Receptor detection rate
$data[0,12] = 1.577167666666667   # Receptor average expression value
$data[0,13] = 4.731503   # Receptor total expression value
$data[0,14] = 0.05487405574265763   # Receptor derived specificity of average expression value
$data[0,15] = 0.05487405574265762   # Receptor derived specificity of total expression value
$data[0,16] = 8.50469579394089   # Edge average expression weight
$data[0,17] = 76.542262145468   # Edge total expression weight
$data[0,18] = 0.002742686339778793   # Edge average expression derived specificity
$data[0,19] = 0.002742686339778793   # Edge total expression derived specificity

# Row 3: ECs -> FAPs
$data[1,0] = "ECs"   # Sending cluster
$data[1,1] = "C1qa"   # Ligand symbol
$data[1,2] = "Cspg4"   # Receptor symbol
$data[1,3] = "FAPs"   # Target cluster
$data[1,4] = 1   # Ligand-expressing cells
$data[1,5] = 0.3333333333333333   # Ligand detection rate
$data[1,6] = 5.392385333333333   # Ligand average expression value
$data[1,7] = 16.177156   # Ligand total expression value
$data[1,8] = 0.04998147672264548   # Ligand derived specificity of average expression value
$data[1,9] = 0.04998147672264548   # Ligand derived specificity of total expression value
$data[1,10] = 3   # Receptor-expressing cells
$data[1,11] = 1   # Receptor detection rate
$data[1,12] = 12.088216   # Receptor average expression value
$data[1,13] = 36.264648   # Receptor total expression value
$data[1,14] = 0.4205827019109694   # Receptor derived specificity of average expression value
$data[1,15] = 0.4205827019109694   # Receptor derived specificity of total expression value
$data[1,16] = 65.18431866456534   # Edge average expression weight
$data[1,17] = 586.6588679810881   # Edge total expression weight
$data[1,18] = 0.02102134452551046   # Edge average expression derived specificity
$data[1,19] = 0.02102134452551046   # Edge total expression derived specificity

# Row 4: ECs -> M2
$data[2,0] = "ECs"   # Sending cluster
$data[2,1] = "C1qa"   # Ligand symbol
$data[2,2] = "Cspg4"   # Receptor symbol
$data[2,3] = "M2"   # Target cluster
$data[2,4] = 1   # Ligand-expressing cells
$data[2,5] = 0.3333333333333333   # Ligand detection rate
$data[2,6] = 5.392385333333333   # Ligand average expression value
$data[2,7] = 16.177156   # Ligand total expression value
$data[2,8] = 0.04998147672264548   # Ligand derived specificity of average expression value
$data[2,9] = 0.04998147672264548   # Ligand derived specificity of total expression value
$data[2,10] = 2   # Receptor-expressing cells
$data[2,11] = 0.6666666666666666   # Receptor detection rate
$data[2,12] = 0.06147399999999999   # Receptor average expression value
$data[2,13] = 0.184422   # Receptor total expression value
$data[2,14] = 0.002138851673172859   # Receptor derived specificity of average expression value
$data[2,15] = 0.002138851673172859   # Receptor derived specificity of total expression value
$data[2,16] = 0.3314914959813333   # Edge average expression weight
$data[2,17] = 2.983423463832   # Edge total expression weight
$data[2,18] = 0.0001069029651158806   # Edge average expression derived specificity
$data[2,19] = 0.0001069029651158806   # Edge total expression derived specificity

# Row 5: ECs -> sCs
$data[3,0] = "ECs"   # Sending cluster
$data[3,1] = "C1qa"   # Ligand symbol
$data[3,2] = "Cspg4"   # Receptor symbol
$data[3,3] = "sCs"   # Target cluster
$data[3,4] = 1   # Ligand-expressing cells
$data[3,5] = 0.3333333333333333   # Ligand detection rate
$data[3,6] = 5.392385333333333   # Ligand average expression value
$data[3,7] = 16.177156   # Ligand total expression value
$data[3,8] = 0.04998147672264548   # Ligand derived specificity of average expression value
$data[3,9] = 0.04998147672264548   # Ligand derived specificity of total expression value
$data[3,10] = 3   # Receptor-expressing cells
$data[3,11] = 1   # Receptor detection rate
$data[3,12] = 15.01473333333333   # Receptor average expression value
$data[3,13] = 45.0442   # Receptor total expression value
$data[3,14] = 0.5224043906732001   # Receptor derived specificity of average expression value
$data[3,15] = 0.5224043906732001   # Receptor derived specificity of total expression value
$data[3,16] = 80.96522781057777   # Edge average expression weight
$data[3,17] = 728.6870502951999   # Edge total expression weight
$data[3,18] = 0.02611054289224035   # Edge average expression derived specificity
$data[3,19] = 0.02611054289224035   # Edge total expression derived specificity

# Row 6: M2 -> ECs
$data[4,0] = "M2"   # Sending cluster
$data[4,1] = "C1qa"   # Ligand symbol
$data[4,2] = "Cspg4"   # Receptor symbol
$data[4,3] = "ECs"   # Target cluster
$data[4,4] = 3   # Ligand-expressing cells
$data[4,5] = 1   # Ligand detection rate
$data[4,6] = 102.49529   # Ligand average expression value
$data[4,7] = 307.48587   # Ligand total expression value
$data[4,8] = 0.9500185232773545   # Ligand derived specificity of average expression value
$data[4,9] = 0.9500185232773545   # Ligand derived specificity of total expression value
$data[4,10] = 3   # Receptor-expressing cells
$data[4,11] = 1   # Receptor detection rate
$data[4,12] = 1.577167666666667   # Receptor average expression value
$data[4,13] = 4.731503   # Receptor total expression value
$data[4,14] = 0.05487405574265763   # Receptor derived specificity of average expression value
$data[4,15] = 0.05487405574265762   # Receptor derived specificity of total expression value
$data[4,16] = 161.6522573736233   # Edge average expression weight
$data[4,17] = 1454.87031636261   # Edge total expression weight
$data[4,18] = 0.05213136940287883   # Edge average expression derived specificity
$data[4,19] = 0.05213136940287882   # Edge total expression derived specificity

# Row 7: M2 -> FAPs
$data[5,0] = "M2"   # Sending cluster
$data[5,1] = "C1qa"   # Ligand symbol
$data[5,2] = "Cspg4"   # Receptor symbol
$data[5,3] = "FAPs"   # Target cluster
$data[5,4] = 3   # Ligand-expressing cells
$data[5,5] = 1   # Ligand detection rate
$data[5,6] = 102.49529   # Ligand average expression value
$data[5,7] = 307.48587   # Ligand total expression value
$data[5,8] = 0.9500185232773545   # Ligand derived specificity of average expression value
$data[5,9] = 0.9500185232773545   # Ligand derived specificity of total expression value
$data[5,10] = 3   # Receptor-expressing cells
$data[5,11] = 1   # Receptor detection rate
$data[5,12] = 12.088216   # Receptor average expression value
$data[5,13] = 36.264648   # Receptor total expression value
$data[5,14] = 0.4205827019109694   # Receptor derived specificity of average expression value
$data[5,15] = 0.4205827019109694   # Receptor derived specificity of total expression value
$data[5,16] = 1238.98520450264   # Edge average expression weight
$data[5,17] = 11150.86684052376   # Edge total expression weight
$data[5,18] = 0.3995613573854589   # Edge average expression derived specificity
$data[5,19] = 0.3995613573854589   # Edge total expression derived specificity

# Row 8: M2 -> M2
$data[6,0] = "M2"   # Sending cluster
$data[6,1] = "C1qa"   # Ligand symbol
$data[6,2] = "Cspg4"   # Receptor symbol
$data[6,3] = "M2"   # Target cluster
$data[6,4] = 3   # Ligand-expressing cells
$data[6,5] = 1   # Ligand detection rate
$data[6,6] = 102.49529   # Ligand average expression value
$data[6,7] = 307.48587   # Ligand total expression value
$data[6,8] = 0.9500185232773545   # Ligand derived specificity of average expression value
$data[6,9] = 0.9500185232773545   # Ligand derived specificity of total expression value
$data[6,10] = 2   # Receptor-expressing cells
$data[6,11] = 0.6666666666666666   # Receptor detection rate
$data[6,12] = 0.06147399999999999   # Receptor average expression value
$data[6,13] = 0.184422   # Receptor total expression value
$data[6,14] = 0.002138851673172859   # Receptor derived specificity of average expression value
$data[6,15] = 0.002138851673172859   # Receptor derived specificity of total expression value
$data[6,16] = 6.30079545746   # Edge average expression weight
$data[6,17] = 56.70715911713999   # Edge total expression weight
$data[6,18] = 0.002031948708056979   # Edge average expression derived specificity
$data[6,19] = 0.002031948708056979   # Edge total expression derived specificity

# Row 9: M2 -> sCs
$data[7,0] = "M2"   # Sending cluster
$data[7,1] = "C1qa"   # Ligand symbol
$data[7,2] = "Cspg4"   # Receptor symbol
$data[7,3] = "sCs"   # Target cluster
$data[7,4] = 3   # Ligand-expressing cells
$data[7,5] = 1   # Ligand detection rate
$data[7,6] = 102.49529   # Ligand average expression value
$data[7,7] = 307.48587   # Ligand total expression value
$data[7,8] = 0.9500185232773545   # Ligand derived specificity of average expression value
$data[7,9] = 0.9500185232773545   # Ligand derived specificity of total expression value
$data[7,10] = 3   # Receptor-expressing cells
$data[7,11] = 1   # Receptor detection rate
$data[7,12] = 15.01473333333333   # Receptor average expression value
$data[7,13] = 45.0442   # Receptor total expression value
$data[7,14] = 0.5224043906732001   # Receptor derived specificity of average expression value
$data[7,15] = 0.5224043906732001   # Receptor derived specificity of total expression value
$data[7,16] = 1538.939447272666   # Edge average expression weight
$data[7,17] = 13850.455025454   # Edge total expression weight
$data[7,18] = 0.4962938477809598   # Edge average expression derived specificity
$data[7,19] = 0.4962938477809598   # Edge total expression derived specificity

$ws.Range("A2:T9").Value = $data

Write-Host "Updated sheet1 with 8 sending/target cluster combinations (rows 2-9)"